$d = $word.ActiveDocument


$xml4 = @'
<w:p w:rsidR="00810BFF" w:rsidRPr="00810BFF" w:rsidRDefault="00810BFF" w:rsidP="00810BFF"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>Joseph Lee</w:t></w:r><w:r w:rsidR="00390E0F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>, g4joseph</w:t></w:r><w:r w:rsidR="00390E0F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:softHyphen/></w:r><w:r w:rsidR="00390E0F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:softHyphen/></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($xml4)

$xml36 = @'
<w:p w:rsidR="00AA6E73" w:rsidRDefault="00A358DA" w:rsidP="00810BFF"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:tab/><w:t xml:space="preserve">We have addressed two security concerns during this assignment: authorization, and preventing </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>DDoS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> attacks. For user authorization and authentication we used Passport, a middleware for Node.js. We used a username and password strategy as our</w:t></w:r><w:r w:rsidR="00AA6E73"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> first form of authentication. The u</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>ser simply enter their username and password into an html form</w:t></w:r><w:r w:rsidR="00AA6E73"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>, their password is then hashed and stored in the database</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>, as storing it in plain-text would result in a catastrophe if a hacker were to ever access the database</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> As a hashed item, it is essentially a one way function that turns the data into a unique string and it cannot be reversed. A single letter from the input returns a completely different hash so it is made safe that way. In addition to hashing, there is added salt to prevent two people from having the same password as it they would both be hashed the same</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> way. With salt, we can have unique hashing where the passwords would be hashed twice by appending addition strings characters.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:softHyphen/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:softHyphen/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$d.Paragraphs(36).Range.InsertXML($xml36)

$xml37 = @'
<w:p w:rsidR="00AA6E73" w:rsidRDefault="00AA6E73" w:rsidP="00810BFF"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">To prevent </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>DDoS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> attacks the ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>ddos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>’ module for Node.js.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Every request made by the same IP address</w:t></w:r><w:r w:rsidR="00A358DA"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>is marked in an internal table using a ‘count’ parameter, like so:</w:t></w:r></w:p>
'@
$d.Paragraphs(37).Range.InsertXML($xml37)

$xml39 = @'
<w:p w:rsidR="00A358DA" w:rsidRDefault="00AA6E73" w:rsidP="00810BFF"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve"> If this count goes above a configurable ‘burst’ number then the expiry parameter doubles. If the count exceeds a pre-determined limit, then the request is denied</w:t></w:r><w:r w:rsidR="00001523"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00001523" w:rsidRPr="00001523"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>The only way for a user who has denied requests to continue is for them to let the expiration time pass, and when expiration hits 0, the entry is deleted from the table, and new requests are allowed like normal.</w:t></w:r></w:p>
'@
$d.Paragraphs(39).Range.InsertXML($xml39)

$xml41 = @'
<w:p w:rsidR="00810BFF" w:rsidRPr="00810BFF" w:rsidRDefault="00810BFF" w:rsidP="00810BFF"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">e. A section that Includes a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>youtube</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t>link ​of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00810BFF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> ​a three minutes video demo of your application (provide username and password if it private). </w:t></w:r></w:p>
'@
$d.Paragraphs(41).Range.InsertXML($xml41)

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
